# "Jos nesto malo dodala" - fill in partner identification + work-package
# descriptions across the partner-budget workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Detaljno budzet
$ws2 = $wb.Worksheets.Item(2)   # Travel - budzet
$ws3 = $wb.Worksheets.Item(3)   # Equipment - budzet
$ws4 = $wb.Worksheets.Item(4)   # Subcontracting - budzet

# --- Fill the short "WPx" labels used on the Travel / Equipment / Subcontracting
#     sheets (rows 5-12, column B) -------------------------------------------------
$wpShort = @("WP1", "WP2", "WP3", "WP4", "WP5", "WP6", "WP7", "WP8")

for ($i = 0; $i -lt $wpShort.Length; $i++) {
    $row = 5 + $i
    $ws2.Range("B$row").Value = $wpShort[$i]
}

for ($i = 0; $i -lt $wpShort.Length; $i++) {
    $row = 5 + $i
    $ws3.Range("B$row").Value = $wpShort[$i]
}

for ($i = 0; $i -lt $wpShort.Length; $i++) {
    $row = 5 + $i
    $ws4.Range("B$row").Value = $wpShort[$i]
}

# --- Replace the generic "WPn-" work-package headings on the main sheet with the
#     full work-package titles -------------------------------------------------
$ws1.Range("A14").Value = "WP1-  Analiza korisničkih zahteva"
$ws1.Range("A15").Value = "WP2- Realizacija neophodnog harvdera"
$ws1.Range("A16").Value = "WP3- Modelovanje sistema"
$ws1.Range("A17").Value = "WP4-  Implementacija veb i mobilne aplikacije"
$ws1.Range("A18").Value = "WP5- Integracija sistema"
$ws1.Range("A19").Value = "WP6- Testiranje"
$ws1.Range("A20").Value = "WP7- Evaluacija i disiminacija"
$ws1.Range("A21").Value = "WP8- Upravljanje projektom"

# --- Partner / project identification -----------------------------------------
$ws1.Range("J4").Value = "Smart Companion"
$ws1.Range("J5").Value = "Elektrotehnički fakultet Univerziteta u Beogradu"
$ws1.Range("J6").Value = "ETF"

# Give column J a bit more room now that it holds real text.
$ws1.Columns.Item(10).ColumnWidth = 16.625

# --- Restore a sensible selection on each sheet, then leave the main budget
#     sheet ("Detaljno budzet") as the active tab ------------------------------
$ws2.Activate()
$ws2.Range("B12").Select()

$ws3.Activate()
$ws3.Range("B5:B12").Select()

$ws4.Activate()
$ws4.Range("B21").Select()

$ws1.Activate()
$ws1.Range("J6").Select()
